# Apply new benchmark measurements and consolidate a duplicate cell style.
#
# Context: the workbook had two nearly-identical fonts/xfs (one keyed to an
# explicit rgb black, the other to the "theme 1 / automatic" black) applied
# to the same visual style (numFmt "#,##0.00", thin border, right-aligned).
# The edit folds the header/data cells that used the "theme" variant onto
# the "rgb" variant's style (D1 already carries that exact xf), and updates
# the measured ratios in columns D/E/F for every benchmark row with new
# multicore/ablation run results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Consolidate the duplicate font/style -----------------------------
# D1 already uses the target xf (numFmt "#,##0.00", border, right align,
# rgb-black font). Copy its format onto every cell that currently uses the
# duplicate "theme black" variant of that same style.
$fmtSource = $ws.Range("D1")
$fmtSource.Copy() | Out-Null
$xlPasteFormats = -4122
$styleFixCells = @(
    "E1","F1","G1",
    "E2","F2","G2",
    "E4","F4","G4",
    "E5","F5","G5",
    "E6","F6",
    "F13","E14"
)
foreach ($ref in $styleFixCells) {
    $ws.Range($ref).PasteSpecial($xlPasteFormats) | Out-Null
}
$excel.CutCopyMode = $false

# --- 2. Update the benchmark ratios (columns D, E, F) ---------------------
$newValues = @{
    "D2"  = 1.133302197;  "E2"  = 1.172297554;  "F2"  = 1.361209652
    "D3"  = 1.362124658;  "E3"  = 1.347774335;  "F3"  = 1.464484041
    "D4"  = 1.057770414;  "E4"  = 1.115796505;  "F4"  = 1.089685778
    "D5"  = 1.159147747;  "E5"  = 1.162994685;  "F5"  = 1.27853834
    "D6"  = 1.121087038;  "E6"  = 1.102806563;  "F6"  = 1.233880283
    "D7"  = 1.294165834;  "E7"  = 1.328733274;  "F7"  = 1.612680826
    "D8"  = 1.319184999;  "E8"  = 1.332107673;  "F8"  = 1.619318667
    "D9"  = 1.392042811;  "E9"  = 1.354171488;  "F9"  = 1.673504762
    "D10" = 1.370564967;  "E10" = 1.311094336;  "F10" = 1.595622405
    "D11" = 1.328475991;  "E11" = 1.333065953;  "F11" = 1.639366909
    "D12" = 1.11130662 ;  "E12" = 1.199347312;  "F12" = 1.401174438
    "D13" = 1.000849953;  "E13" = 1.007662151;  "F13" = 1.013618789
    "D14" = 1.312139042;  "E14" = 1.28302802 ;  "F14" = 1.542333722
    "D15" = 1.121908297;  "E15" = 1.121553102;  "F15" = 1.269679789
    "D16" = 1.220290755;  "E16" = 1.226602354;  "F16" = 1.410047258
}

foreach ($ref in $newValues.Keys) {
    $ws.Range($ref).Value2 = $newValues[$ref]
}
